# ============================================================================
# Edit: add 2022-Q1 quarterly holdings sheet, update 总计 (Total) summary sheet
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Step 1: the existing "总计" sheet becomes the new "2022-Q1" data sheet.
# This lets it keep sheetId=4 (as if a sheet had been freshly inserted right
# before "总计") while a brand-new sheet (sheetId=5) takes over the "总计"
# name/position - matching the sheet-id layout of the target workbook.
# ----------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# ----------------------------------------------------------------------------
# Step 2: clear whatever is currently on the (repurposed) "2022-Q1" sheet and
# lay out the fund holdings table: 基金代码/基金名称/基金规模/股票总仓位/
# 仓位占比/持有市值(亿元)/仓位排名
# ----------------------------------------------------------------------------
$q1.Cells.Clear()

$styleSrc = $wb.Worksheets.Item("2021-Q4")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q1.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$fundRows = @(
    @(0, '159920', '华夏恒生ETF(QDII)', '151.31', '95.19', '7.92', '11.9838', 1),
    @(1, '513550', '华泰柏瑞中证港股通50ETF', '31.28', '98.89', '10.39', '3.2500', 1),
    @(2, '513660', '华夏沪港通恒生ETF', '19.61', '97.34', '8.49', '1.6649', 2),
    @(3, '501025', '鹏华港股通中证香港银行投资指数（LOF）A', '9.81', '94.47', '15.53', '1.5235', 1),
    @(4, '010365', '鹏华港股通中证香港银行投资指数（LOF）C', '6.07', '94.47', '15.53', '0.9427', 1),
    @(5, '007139', '富国民裕进取沪港深成长精选混合', '12.79', '92.21', '4.71', '0.6024', 9),
    @(6, '513600', '南方恒生ETF', '5.89', '99.00', '8.48', '0.4995', 3),
    @(7, '517300', '国寿安保中证沪港深300ETF', '25.13', '99.04', '1.84', '0.4624', 6),
    @(8, '006809', '泰康港股通中证香港银行投资指数A', '1.99', '94.73', '15.53', '0.3090', 1),
    @(9, '007354', '创金合信港股通量化股票A', '3.84', '91.20', '5.92', '0.2273', 2),
    @(10, '513900', '华安CES港股通精选100ETF', '2.15', '96.24', '10.10', '0.2172', 1),
    @(11, '010204', '中银港股通优势成长股票', '3.19', '83.00', '5.51', '0.1758', 3),
    @(12, '517000', '银华中证沪港深500ETF', '6.14', '94.65', '2.58', '0.1584', 3),
    @(13, '517080', '汇添富中证沪港深500ETF', '5.68', '91.59', '2.57', '0.1460', 3),
    @(14, '006810', '泰康港股通中证香港银行投资指数C', '0.90', '94.73', '15.53', '0.1398', 1),
    @(15, '517100', '富国中证沪港深500ETF', '4.13', '99.22', '2.71', '0.1119', 3),
    @(16, '159712', '国泰中证港股通50ETF', '0.85', '95.21', '9.45', '0.0803', 1),
    @(17, '513990', '招商上证港股通ETF', '0.59', '96.48', '7.67', '0.0453', 2),
    @(18, '005142', '中融沪港深大消费主题灵活配置混合A', '0.49', '88.98', '5.27', '0.0258', 6),
    @(19, '501309', '国泰恒生港股通指数（LOF）', '0.36', '92.35', '6.79', '0.0244', 2),
    @(20, '162416', '华宝港股通恒生香港35指数(LOF)', '0.21', '94.50', '9.68', '0.0203', 2),
    @(21, '005143', '中融沪港深大消费主题灵活配置混合C', '0.33', '88.98', '5.27', '0.0174', 6),
    @(22, '160925', '大成中华沪深港300指数（LOF）A', '0.54', '93.14', '2.97', '0.0160', 3),
    @(23, '007357', '创金合信港股通量化股票C', '0.26', '91.20', '5.92', '0.0154', 2),
    @(24, '005701', '上投摩根香港精选港股通混合', '0.44', '84.37', '3.06', '0.0135', 10),
    @(25, '166402', '浦银安盛中证锐联沪港深基本面100指数（LOF）', '0.19', '90.95', '6.55', '0.0124', 1),
    @(26, '517010', '易方达中证沪港深500交易型开放式指数证券投资基金', '0.44', '91.01', '2.72', '0.0120', 3),
    @(27, '006106', '景顺长城量化港股通股票', '0.34', '85.20', '3.43', '0.0117', 3),
    @(28, '005707', '富国港股通量化精选股票', '0.24', '80.43', '4.76', '0.0114', 2),
    @(29, '517170', '华夏中证沪港深500交易型开放式指数证券投资基金', '0.33', '94.49', '2.61', '0.0086', 3),
    @(30, '008973', '大成中华沪深港300指数(LOF)C', '0.02', '93.14', '2.97', '0.0006', 3)
)

foreach ($row in $fundRows) {
    $r = [int]$row[0] + 2
    $q1.Cells.Item($r, 1).Value = [int]$row[0]
    $q1.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1.Cells.Item($r, 6).Value = "'" + $row[5]
    $q1.Cells.Item($r, 7).Value = "'" + $row[6]
    $q1.Cells.Item($r, 8).Value = [int]$row[7]
}

# Re-apply the same formatting (bold / centered / thin-bordered header+index
# style) that the other quarterly sheets use, by broadcasting copied formats
# from a single already-styled cell onto the destination ranges.
$styleSrc.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$q1.Range("A2:A32").PasteSpecial(-4122)

$q1.Range("A1").Select()

# ----------------------------------------------------------------------------
# Step 3: populate the new "总计" sheet - same 日期/持有数量(只)/持有市值(亿元)
# summary as before, with a new row for 2022-Q1 inserted at the top.
# ----------------------------------------------------------------------------
$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$summaryRows = @(
    @(0, "2022-Q1", 31, 22.73),
    @(1, "2021-Q4", 26, 15.67),
    @(2, "2021-Q3", 21, 10.8),
    @(3, "2021-Q2", 27, 16.4)
)

foreach ($row in $summaryRows) {
    $r = [int]$row[0] + 2
    $total.Cells.Item($r, 1).Value = [int]$row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = [int]$row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
}

$styleSrc.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)

$total.Range("A1").Select()
